$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "About": insert three new rows (two notes + one blank spacer) right
# before the existing "See ""cpi.xlsx""..." footnote block, pushing the rest
# of the footnote rows down by three, and add the new explanatory note.
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Rows("18:20").Insert()

$wsAbout.Range("A18").Value = "The model uses LDVs elasticity for all vehicle types because no data on price elasticity"
$wsAbout.Range("A19").Value = "of other vehicle types with respect to fuel economy is available."
$wsAbout.Range("A18:A19").Font.Bold = $false

# ---------------------------------------------------------------------------
# Sheet "EoVPwFE": relabel the elasticity cell, make it wrap, and resize the
# header / data rows to fit the new two-line label.
# ---------------------------------------------------------------------------
$wsEoV = $wb.Worksheets.Item("EoVPwFE")

$wsEoV.Range("B1").Value = "Elasticity (dimensionless)"
$wsEoV.Range("B1").WrapText = $true
$wsEoV.Range("B1").HorizontalAlignment = -4152
$wsEoV.Rows(1).RowHeight = 45
$wsEoV.Rows(2).RowHeight = 14.45

# ---------------------------------------------------------------------------
# Restore the on-screen selections so "About" is the active/selected sheet
# (with the newly-inserted block selected) and "EoVPwFE" remembers B1 as its
# last selection.
# ---------------------------------------------------------------------------
$wsEoV.Activate() | Out-Null
$wsEoV.Range("B1").Select() | Out-Null

$wsAbout.Activate() | Out-Null
$wsAbout.Range("A20:XFD23").Select() | Out-Null
